$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("YDS")

$ws.Range("B2").Value = "NIG(0.8251779309585252, 0.6012245913602724, 1.2612175380261814, 3.1332032875500233)"
$ws.Range("C2").Value = "JSU(-0.8857409098576617, 1.0557264788051692, 4.3742083693798905, 4.4378454690354125)"
$ws.Range("D2").Value = "NIG(0.5423499285266689, 0.28352982601061405, 2.801250485877909, 2.8251319740907905)"
$ws.Range("E2").Value = "JSU(-1.0807344524313427, 1.1025791924412918, 3.248318519970014, 4.821470921730919)"
